# Insert a new data row for "Terminal Hortofrutícola Agro Chillán" / Mango
# at worksheet row 71, pushing the existing rows 71-138 down to 72-139.
# (Net effect observed in the target diff: dimension grows from A1:T138 to
#  A1:T139, and every row from 71 onward shifts down by one, with the new
#  row containing a fresh weekly observation.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 71..138 down to 72..139, creating a blank row 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new observation.
$ws.Range("A71").Value = 7
$ws.Range("B71").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C71").Value = "Ñuble"
$ws.Range("D71").Value = 45079
$ws.Range("D71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100108
$ws.Range("H71").Value = "Tropicales y subtropicales"
$ws.Range("I71").Value = 100108002
$ws.Range("J71").Value = "Mango"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 50
$ws.Range("N71").Value = 9000
$ws.Range("O71").Value = 9000
$ws.Range("P71").Value = 9000
$ws.Range("Q71").Value = "$/bandeja 4 kilos"
$ws.Range("R71").Value = "Perú"
$ws.Range("S71").Value = 2250
$ws.Range("T71").Value = 4
